# Regression Suite for CDS
# TC01_CDS_phs001713.xlsx - "ParticipantsTab" (row 2 / cell B2) query gains a
# trailing "ORDER BY participant_id desc" clause before the final LIMIT, and the
# saved view/selection now lands on B2 instead of B4.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$participantsQuery = @'
WITH Distinct_Samples AS (
    SELECT DISTINCT
        sp.participant_id,
        sp.study_participant_id,
        s.study_name,
        s.phs_accession,
        sp.gender,
        smp.sample_id
    FROM 
        df_participant sp
    JOIN 
        df_study s ON sp."study.phs_accession" = s.phs_accession
    JOIN 
        df_sample smp ON smp."participant.study_participant_id" = sp.study_participant_id
    JOIN
        df_diagnosis d ON d."participant.study_participant_id" = sp.study_participant_id
    JOIN
        df_program p ON p.program_acronym = s."program.program_acronym"
    JOIN
        df_file f1 ON f1."sample.sample_id" = smp.sample_id
    JOIN
        df_genomic_info gi ON gi."file.file_id" = f1.file_id
    WHERE 
        s.phs_accession = 'phs001713'
),
Sample_Limit AS (
    SELECT
        participant_id,
        study_participant_id,
        study_name,
        phs_accession,
        gender,
        sample_id,
        ROW_NUMBER() OVER (PARTITION BY study_participant_id ORDER BY sample_id) as rn
    FROM 
        Distinct_Samples
)
SELECT
    participant_id AS "Participant ID", 
    study_name AS "Study Name",
    phs_accession AS Accession,
    gender AS Gender,
    GROUP_CONCAT(
        CASE 
            WHEN rn <= 5 THEN sample_id 
        END, ', '
    ) ||
    CASE 
        WHEN MAX(rn) > 5 THEN ', ...' 
        ELSE '' 
    END AS Samples
FROM 
    Sample_Limit
GROUP BY
    participant_id, 
    study_name,
    phs_accession,
    gender
ORDER BY participant_id desc
LIMIT 100;
'@

$ws.Range("B2").Value = $participantsQuery

# The cell already wraps text at style "3"; setting the (longer) value makes the
# native autofit grow past Excel's 409.5pt row-height ceiling, so pin it back to
# the sheet's existing max row height (matches rows 3 and 4).
$ws.Rows.Item(2).RowHeight = 409.5

# Move the saved selection / scroll position from B4 to B2.
$ws.Range("B2").Select()

